# "Stylesheet and T5 work"
# Corrects the Transcription numbers attributed to Hand9/Hand10 on the
# "Hand " sheet, adds three new hands (16-18) with their attributions,
# removes the stray "title" note on the "ana" sheet and appends a new
# "prep" part-of-speech entry there.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "Hand " sheet (transcription attributions + new hands)
# ------------------------------------------------------------------
$hand = $wb.Worksheets.Item("Hand ")

# Swap which transcription each of Hand9 / Hand10 belongs to, and give
# Hand11 (Domhnall MacMhuirich) the same Transcription 9 reference.
$hand.Cells.Item(9, 3).Value = "Transcription 12"
$hand.Cells.Item(10, 3).Value = "Transcription 9"
$hand.Cells.Item(11, 3).Value = "Transcription 9"

# New hands 16-18 - fill column by column (Hand ids, then names, then
# transcriptions) to match how the sheet was actually populated.
$hand.Cells.Item(16, 1).Value = "Hand16"
$hand.Cells.Item(17, 1).Value = "Hand17"
$hand.Cells.Item(18, 1).Value = "Hand18"

$hand.Cells.Item(16, 2).Value = "Rev. John Beaton"
$hand.Cells.Item(17, 2).Value = "Neil Beaton"
$hand.Cells.Item(18, 2).Value = "Dubhghall Albanach mac mhic Cathail"

$hand.Cells.Item(16, 3).Value = "Transcription 5"
$hand.Cells.Item(17, 3).Value = "Transcription 5"
$hand.Cells.Item(18, 3).Value = "Transcription 4"

$hand.Range("C11").Select() | Out-Null

# ------------------------------------------------------------------
# "ana" sheet (word-class list)
# ------------------------------------------------------------------
$ana = $wb.Worksheets.Item("ana")

# "title" was mistakenly attached next to "tit" - remove it.
$ana.Cells.Item(17, 2).ClearContents()

# Add the missing "prep" entry.
$ana.Cells.Item(20, 1).Value = "prep"

$ana.Range("B17").Select() | Out-Null

# ------------------------------------------------------------------
# Re-select the "Glyphs" sheet (left as the active tab) at A9.
# ------------------------------------------------------------------
$glyphs = $wb.Worksheets.Item("Glyphs")
$glyphs.Activate() | Out-Null
$glyphs.Range("A9").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 19
} catch {
}
